$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.870.24'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -5.40%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.113.29'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.93%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.06%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.42'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -10.00%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.586'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -10.50%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.111.88'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.91%  '

# Row 10
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.89%  '

# Row 11
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.116'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.47%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.378'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.17%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.656.81'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.95%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.127'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.92%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.048.71'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.15%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '24.53'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.99%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.118.16'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.49%  '

# Row 18
$ws.Range('E18').Value = '  -6.28%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '404.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.47%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.38%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.72%  '

# Row 22
$ws.Range('E22').Value = '  -3.96%  '

# Row 23
$ws.Range('E23').Value = '  -0.08%  '

# Row 24
$ws.Range('E24').Value = '  +0.22%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '68.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.05%  '

# Row 26
$ws.Range('E26').Value = '  -2.54%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.488'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.27%  '

# Row 28
$ws.Range('E28').Value = '  -11.57%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.68'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.64%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.03%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.26'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.10%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.78'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.19%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.86'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.77%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.21'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.06%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '154.62'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.43%  '

# Row 37
$ws.Range('E37').Value = '  -7.09%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.56%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.726.57'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.50%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.64'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.72%  '

# Row 41
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.62'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -10.29%  '

# Row 42
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.10'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.29%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '38.17'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.86%  '

# Row 44
$ws.Range('E44').Value = '  -7.62%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0605'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.30%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0255'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.15%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.23'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -11.36%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '285.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -8.49%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.88'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -9.60%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.04%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0969'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.82%  '
